$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.595.64"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "3.251.21"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Formula = "'580.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Formula = "'184.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Formula = "'0.606"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "3.251.03"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").Formula = "'6.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("D13").Value = "3.799.62"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("D16").Value = "67.619.86"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "3.214.09"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").Formula = "'13.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Formula = "'394.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.23%  "
$ws.Range("D22").Formula = "'7.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D24").Formula = "'71.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Formula = "'0.517"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Formula = "'0.186"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("D28").Formula = "'9.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Formula = "'5.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.26%  "
$ws.Range("D32").Formula = "'22.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Formula = "'161.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Formula = "'1.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").Formula = "'26.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("D43").Formula = "'2.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.97%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Formula = "'40.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").Value = "2.615.48"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").Formula = "'24.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("D48").Formula = "'334.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").Formula = "'6.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("E51").Value = "  -0.45%  "
